# Update currentAveragePrice / LevePrice / LeveProfit figures (cols H-N)
# across the Leve-profit tables on several sheets, per the scheduled
# market-board data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 112.875
$ws.Range("I42").Value = 13.5
$ws.Range("J42").Value = 212.25
$ws.Range("K42").Value = 40.5
$ws.Range("L42").Value = 636.75
$ws.Range("M42").Value = 189.5
$ws.Range("N42").Value = -1096.75
$ws.Range("H52").Value = 3150
$ws.Range("I52").Value = 1000
$ws.Range("J52").Value = 3418.75
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 10256.25
$ws.Range("M52").Value = -2840
$ws.Range("N52").Value = -10576.25
$ws.Range("H112").Value = 1112.6305
$ws.Range("J112").Value = 1124.525
$ws.Range("L112").Value = 3373.575
$ws.Range("N112").Value = -5589.575000000001
$ws.Range("H141").Value = 6608.1113
$ws.Range("I141").Value = 2229.8
$ws.Range("J141").Value = 28499.666
$ws.Range("K141").Value = 6689.400000000001
$ws.Range("L141").Value = 85498.99800000001
$ws.Range("M141").Value = -1509.400000000001
$ws.Range("N141").Value = -95858.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1201.0869
$ws.Range("I2").Value = 888.6667
$ws.Range("J2").Value = 1786.875
$ws.Range("K2").Value = 888.6667
$ws.Range("L2").Value = 1786.875
$ws.Range("M2").Value = -775.6667
$ws.Range("N2").Value = -2012.875
$ws.Range("H61").Value = 5419.229
$ws.Range("I61").Value = 6874.6665
$ws.Range("J61").Value = 3963.7917
$ws.Range("K61").Value = 6874.6665
$ws.Range("L61").Value = 3963.7917
$ws.Range("M61").Value = -6662.6665
$ws.Range("N61").Value = -4387.7917
$ws.Range("H74").Value = 1358.8611
$ws.Range("I74").Value = 1350.931
$ws.Range("J74").Value = 1391.7142
$ws.Range("K74").Value = 1350.931
$ws.Range("L74").Value = 1391.7142
$ws.Range("M74").Value = -476.931
$ws.Range("N74").Value = -3139.7142
$ws.Range("H77").Value = 1358.8611
$ws.Range("I77").Value = 1350.931
$ws.Range("J77").Value = 1391.7142
$ws.Range("K77").Value = 6754.655000000001
$ws.Range("L77").Value = 6958.571
$ws.Range("M77").Value = -2386.655000000001
$ws.Range("N77").Value = -15694.571
$ws.Range("H116").Value = 1201.0869
$ws.Range("I116").Value = 888.6667
$ws.Range("J116").Value = 1786.875
$ws.Range("K116").Value = 888.6667
$ws.Range("L116").Value = 1786.875
$ws.Range("M116").Value = 1405.3333
$ws.Range("N116").Value = -6374.875
$ws.Range("H136").Value = 5419.229
$ws.Range("I136").Value = 6874.6665
$ws.Range("J136").Value = 3963.7917
$ws.Range("K136").Value = 20623.9995
$ws.Range("L136").Value = 11891.3751
$ws.Range("M136").Value = -18073.9995
$ws.Range("N136").Value = -16991.3751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1201.0869
$ws.Range("I3").Value = 888.6667
$ws.Range("J3").Value = 1786.875
$ws.Range("K3").Value = 888.6667
$ws.Range("L3").Value = 1786.875
$ws.Range("M3").Value = -774.6667
$ws.Range("N3").Value = -2014.875
$ws.Range("H64").Value = 511.4375
$ws.Range("I64").Value = 506.9091
$ws.Range("J64").Value = 521.4
$ws.Range("K64").Value = 506.9091
$ws.Range("L64").Value = 521.4
$ws.Range("M64").Value = -281.9091
$ws.Range("N64").Value = -971.4
$ws.Range("H67").Value = 511.4375
$ws.Range("I67").Value = 506.9091
$ws.Range("J67").Value = 521.4
$ws.Range("K67").Value = 506.9091
$ws.Range("L67").Value = 521.4
$ws.Range("M67").Value = 273.0909
$ws.Range("N67").Value = -2081.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4413.4595
$ws.Range("I134").Value = 4552.8
$ws.Range("J134").Value = 1975
$ws.Range("K134").Value = 13658.4
$ws.Range("L134").Value = 5925
$ws.Range("M134").Value = -11123.4
$ws.Range("N134").Value = -10995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2300
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H113").Value = 1430
$ws.Range("I113").Value = 1705.1333
$ws.Range("J113").Value = 1223.65
$ws.Range("K113").Value = 5115.3999
$ws.Range("L113").Value = 3670.95
$ws.Range("M113").Value = -2945.3999
$ws.Range("N113").Value = -8010.950000000001
$ws.Range("H136").Value = 39214.477
$ws.Range("I136").Value = 2409.0715
$ws.Range("J136").Value = 112825.29
$ws.Range("K136").Value = 7227.2145
$ws.Range("L136").Value = 338475.87
$ws.Range("M136").Value = -2127.2145
$ws.Range("N136").Value = -348675.87
$ws.Range("H141").Value = 1812.5
$ws.Range("J141").Value = 2000
$ws.Range("L141").Value = 6000
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5042.909
$ws.Range("I80").Value = 5060.3125
$ws.Range("J80").Value = 4996.5
$ws.Range("K80").Value = 5060.3125
$ws.Range("L80").Value = 4996.5
$ws.Range("M80").Value = -4062.3125
$ws.Range("N80").Value = -6992.5
$ws.Range("H83").Value = 5042.909
$ws.Range("I83").Value = 5060.3125
$ws.Range("J83").Value = 4996.5
$ws.Range("K83").Value = 25301.5625
$ws.Range("L83").Value = 24982.5
$ws.Range("M83").Value = -20309.5625
$ws.Range("N83").Value = -34966.5
$ws.Range("H97").Value = 2203.3
$ws.Range("I97").Value = 2196.1538
$ws.Range("J97").Value = 2216.5715
$ws.Range("K97").Value = 2196.1538
$ws.Range("L97").Value = 2216.5715
$ws.Range("M97").Value = -1700.1538
$ws.Range("N97").Value = -3208.5715
$ws.Range("H126").Value = 2258.2144
$ws.Range("I126").Value = 1885.8
$ws.Range("J126").Value = 3189.25
$ws.Range("K126").Value = 5657.4
$ws.Range("L126").Value = 9567.75
$ws.Range("M126").Value = -3187.4
$ws.Range("N126").Value = -14507.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 18521002
$ws.Range("I122").Value = 22224602
$ws.Range("K122").Value = 66673806
$ws.Range("M122").Value = -66671356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 41622.08
$ws.Range("I122").Value = 56986.5
$ws.Range("J122").Value = 2113.5715
$ws.Range("K122").Value = 170959.5
$ws.Range("L122").Value = 6340.7145
$ws.Range("M122").Value = -168509.5
$ws.Range("N122").Value = -11240.7145
$ws.Range("H135").Value = 72432.89
$ws.Range("J135").Value = 72432.89
$ws.Range("L135").Value = 72432.89
$ws.Range("N135").Value = -82572.89
$ws.Range("H136").Value = 6806561
$ws.Range("I136").Value = 30303806
$ws.Range("J136").Value = 4727.079
$ws.Range("K136").Value = 90911418
$ws.Range("L136").Value = 14181.237
$ws.Range("M136").Value = -90908868
$ws.Range("N136").Value = -19281.237
